$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (shared string cell A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 14:42"

# Update country statistics rows.
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2389166
$ws.Range("C4").Value = 1013
$ws.Range("D4").Value = 1003258
$ws.Range("E4").Value = 1263274
$ws.Range("G4").Value = 24
$ws.Range("H4").Value = 122634

# Row 5 - Brasil
$ws.Range("B5").Value = 1113606
$ws.Range("C5").Value = 2258
$ws.Range("E5").Value = 468095

# Row 14 - Alemania
$ws.Range("D14").Value = 175700
$ws.Range("E14").Value = 7450

# Row 31 - Paises Bajos
$ws.Range("B31").Value = 49722
$ws.Range("C31").Value = 64
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 6095

# Row 33 - Emiratos Arabes Unidos
$ws.Range("B33").Value = 45683
$ws.Range("C33").Value = 380
$ws.Range("D33").Value = 33703
$ws.Range("E33").Value = 11675
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 305

# Row 37 - Portugal
$ws.Range("B37").Value = 39737
$ws.Range("C37").Value = 345
$ws.Range("D37").Value = 25829
$ws.Range("E37").Value = 12368
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 1540

# Row 63 - Dinamarca
$ws.Range("B63").Value = 12561
$ws.Range("C63").Value = 34
$ws.Range("D63").Value = 11393
$ws.Range("E63").Value = 565
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 603

# Row 76 - Uzbekistan
$ws.Range("B76").Value = 6535
$ws.Range("C76").Value = 74
$ws.Range("D76").Value = 4520
$ws.Range("E76").Value = 1996

# Row 82 - Guinea
$ws.Range("B82").Value = 5040
$ws.Range("C82").Value = 52
$ws.Range("D82").Value = 3685
$ws.Range("E82").Value = 1327
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 28

# Row 92 - Bosnia y Herzegovina
$ws.Range("B92").Value = 3588
$ws.Range("C92").Value = 63
$ws.Range("D92").Value = 2285
$ws.Range("E92").Value = 1131
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 172

# Row 101 - Croacia
$ws.Range("B101").Value = 2366
$ws.Range("C101").Value = 30
$ws.Range("E101").Value = 117
